$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 32000
$ws.Range("J87").Value = 32000
$ws.Range("L87").Value = 32000
$ws.Range("N87").Value = -34496
$ws.Range("H90").Value = 32000
$ws.Range("J90").Value = 32000
$ws.Range("L90").Value = 96000
$ws.Range("N90").Value = -108480
$ws.Range("H92").Value = 3654.4092
$ws.Range("I92").Value = 1464.85
$ws.Range("J92").Value = 25550
$ws.Range("K92").Value = 1464.85
$ws.Range("L92").Value = 25550
$ws.Range("M92").Value = -216.8499999999999
$ws.Range("N92").Value = -28046
$ws.Range("H125").Value = 1434.7368
$ws.Range("I125").Value = 2376.4
$ws.Range("J125").Value = 388.44446
$ws.Range("K125").Value = 21387.6
$ws.Range("L125").Value = 3496.00014
$ws.Range("M125").Value = -18927.6
$ws.Range("N125").Value = -8416.00014
$ws.Range("H132").Value = 4723165
$ws.Range("I132").Value = 5671.9062
$ws.Range("J132").Value = 11911726
$ws.Range("K132").Value = 17015.7186
$ws.Range("L132").Value = 35735178
$ws.Range("M132").Value = -14485.7186
$ws.Range("N132").Value = -35740238
$ws.Range("H137").Value = 5002728.5
$ws.Range("I137").Value = 1629.7142
$ws.Range("J137").Value = 10530259
$ws.Range("K137").Value = 4889.142599999999
$ws.Range("L137").Value = 31590777
$ws.Range("M137").Value = -2339.142599999999
$ws.Range("N137").Value = -31595877

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1400
$ws.Range("I2").Value = 1266.6666
$ws.Range("J2").Value = 1514.2858
$ws.Range("K2").Value = 1266.6666
$ws.Range("L2").Value = 1514.2858
$ws.Range("M2").Value = -1153.6666
$ws.Range("N2").Value = -1740.2858
$ws.Range("H32").Value = 7694.0845
$ws.Range("I32").Value = 8212.375
$ws.Range("K32").Value = 8212.375
$ws.Range("M32").Value = -7925.375
$ws.Range("H74").Value = 6667829
$ws.Range("I74").Value = 7693109
$ws.Range("J74").Value = 3509.8
$ws.Range("K74").Value = 7693109
$ws.Range("L74").Value = 3509.8
$ws.Range("M74").Value = -7692235
$ws.Range("N74").Value = -5257.8
$ws.Range("H77").Value = 6667829
$ws.Range("I77").Value = 7693109
$ws.Range("J77").Value = 3509.8
$ws.Range("K77").Value = 38465545
$ws.Range("L77").Value = 17549
$ws.Range("M77").Value = -38461177
$ws.Range("N77").Value = -26285
$ws.Range("H116").Value = 1400
$ws.Range("I116").Value = 1266.6666
$ws.Range("J116").Value = 1514.2858
$ws.Range("K116").Value = 1266.6666
$ws.Range("L116").Value = 1514.2858
$ws.Range("M116").Value = 1027.3334
$ws.Range("N116").Value = -6102.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1400
$ws.Range("I3").Value = 1266.6666
$ws.Range("J3").Value = 1514.2858
$ws.Range("K3").Value = 1266.6666
$ws.Range("L3").Value = 1514.2858
$ws.Range("M3").Value = -1152.6666
$ws.Range("N3").Value = -1742.2858
$ws.Range("I80").Value = 18408
$ws.Range("J80").Value = 379.33334
$ws.Range("K80").Value = 18408
$ws.Range("L80").Value = 379.33334
$ws.Range("M80").Value = -17410
$ws.Range("N80").Value = -2375.33334
$ws.Range("I83").Value = 18408
$ws.Range("J83").Value = 379.33334
$ws.Range("K83").Value = 92040
$ws.Range("L83").Value = 1896.6667
$ws.Range("M83").Value = -87048
$ws.Range("N83").Value = -11880.6667
$ws.Range("H86").Value = 19232590
$ws.Range("J86").Value = 125001820
$ws.Range("L86").Value = 125001820
$ws.Range("N86").Value = -125004066
$ws.Range("H89").Value = 19232590
$ws.Range("J89").Value = 125001820
$ws.Range("L89").Value = 625009100
$ws.Range("N89").Value = -625020332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2294.0881
$ws.Range("I58").Value = 1303.8462
$ws.Range("J58").Value = 2907.0952
$ws.Range("K58").Value = 1303.8462
$ws.Range("L58").Value = 2907.0952
$ws.Range("M58").Value = -1100.8462
$ws.Range("N58").Value = -3313.0952
$ws.Range("H99").Value = 1234.0476
$ws.Range("I99").Value = 1154.6923
$ws.Range("J99").Value = 1363
$ws.Range("K99").Value = 1154.6923
$ws.Range("L99").Value = 1363
$ws.Range("M99").Value = 343.3077000000001
$ws.Range("N99").Value = -4359
$ws.Range("H126").Value = 1234.0476
$ws.Range("I126").Value = 1154.6923
$ws.Range("J126").Value = 1363
$ws.Range("K126").Value = 3464.0769
$ws.Range("L126").Value = 4089
$ws.Range("M126").Value = -994.0769
$ws.Range("N126").Value = -9029
$ws.Range("H136").Value = 2294.0881
$ws.Range("I136").Value = 1303.8462
$ws.Range("J136").Value = 2907.0952
$ws.Range("K136").Value = 3911.5386
$ws.Range("L136").Value = 8721.285600000001
$ws.Range("M136").Value = -1361.5386
$ws.Range("N136").Value = -13821.2856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1157558.1
$ws.Range("I2").Value = 248
$ws.Range("J2").Value = 1984208.2
$ws.Range("K2").Value = 1488
$ws.Range("L2").Value = 11905249.2
$ws.Range("M2").Value = -1375
$ws.Range("N2").Value = -11905475.2
$ws.Range("H38").Value = 150.2
$ws.Range("I38").Value = 802.5
$ws.Range("J38").Value = 49.846153
$ws.Range("K38").Value = 2407.5
$ws.Range("L38").Value = 149.538459
$ws.Range("M38").Value = -2060.5
$ws.Range("N38").Value = -843.538459
$ws.Range("H82").Value = 4942.3335
$ws.Range("H85").Value = 4942.3335
$ws.Range("H86").Value = 1155.9166
$ws.Range("I86").Value = 924
$ws.Range("J86").Value = 1321.5714
$ws.Range("K86").Value = 2772
$ws.Range("L86").Value = 3964.7142
$ws.Range("M86").Value = -1586
$ws.Range("N86").Value = -6336.7142
$ws.Range("H89").Value = 1155.9166
$ws.Range("I89").Value = 924
$ws.Range("J89").Value = 1321.5714
$ws.Range("K89").Value = 8316
$ws.Range("L89").Value = 11894.1426
$ws.Range("M89").Value = -2388
$ws.Range("N89").Value = -23750.1426
$ws.Range("H140").Value = 3929.037
$ws.Range("I140").Value = 1765.6666
$ws.Range("K140").Value = 5296.9998
$ws.Range("M140").Value = -116.9997999999996
$ws.Range("H141").Value = 8736.632
$ws.Range("I141").Value = 8775.714
$ws.Range("J141").Value = 8713.833000000001
$ws.Range("K141").Value = 26327.142
$ws.Range("L141").Value = 26141.499
$ws.Range("M141").Value = -21147.142
$ws.Range("N141").Value = -36501.499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3293.7812
$ws.Range("I126").Value = 1641.3529
$ws.Range("J126").Value = 5166.533
$ws.Range("K126").Value = 4924.0587
$ws.Range("L126").Value = 15499.599
$ws.Range("M126").Value = -2454.0587
$ws.Range("N126").Value = -20439.599

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6566.875
$ws.Range("I7").Value = 7800
$ws.Range("J7").Value = 5686.0713
$ws.Range("K7").Value = 7800
$ws.Range("L7").Value = 5686.0713
$ws.Range("M7").Value = -7688
$ws.Range("N7").Value = -5910.0713
$ws.Range("H40").Value = 4193.3613
$ws.Range("I40").Value = 4648.85
$ws.Range("J40").Value = 3624
$ws.Range("K40").Value = 4648.85
$ws.Range("L40").Value = 3624
$ws.Range("M40").Value = -4512.85
$ws.Range("N40").Value = -3896
$ws.Range("H126").Value = 6566.875
$ws.Range("I126").Value = 7800
$ws.Range("J126").Value = 5686.0713
$ws.Range("K126").Value = 23400
$ws.Range("L126").Value = 17058.2139
$ws.Range("M126").Value = -20930
$ws.Range("N126").Value = -21998.2139
$ws.Range("H132").Value = 9811163
$ws.Range("I132").Value = 5235.5356
$ws.Range("J132").Value = 21748814
$ws.Range("K132").Value = 15706.6068
$ws.Range("L132").Value = 65246442
$ws.Range("M132").Value = -13176.6068
$ws.Range("N132").Value = -65251502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1087.3438
$ws.Range("I136").Value = 963.8461
$ws.Range("J136").Value = 1622.5
$ws.Range("K136").Value = 2891.5383
$ws.Range("L136").Value = 4867.5
$ws.Range("M136").Value = -341.5383000000002
$ws.Range("N136").Value = -9967.5
$ws.Range("H138").Value = 60417.668
$ws.Range("J138").Value = 60417.668
$ws.Range("L138").Value = 60417.668
$ws.Range("N138").Value = -70697.66800000001
